$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F
$ws.Range("F1").Value = "correctAnswer"

# correctAnswer values for rows 2-9 (fixed expected response, no counterbalancing yet)
$answers = @(1, 2, 1, 2, 1, 1, 2, 2)
for ($i = 0; $i -lt $answers.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 6).Value = $answers[$i]
}

# Move selection to match the post-edit active cell
$ws.Range("F11").Select()
